$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old shopee links in A2:A27 but keep each cell's existing style.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
}

# A3 picks up the plain (border/fill-less) "Hyperlink" cell style even
# though it ends up with no link of its own.
$ws.Range("A3").Style = "Hyperlink"

# A5 becomes the one remaining link - a Shopee product title that is wired
# up as a real hyperlink to the product page. Apply the Hyperlink style
# first so Hyperlinks.Add doesn't re-introduce the border/fill look.
$ws.Range("A5").Style = "Hyperlink"

$title = "combo 10 quần lót nữ su đúc cao cấp cạp cao tàn hình chữ GODDESS siêu đẹp | Shopee Việt Nam"
$ws.Range("A5").Value = $title

$url = "https://shopee.vn/combo-10-qu%E1%BA%A7n-l%C3%B3t-n%E1%BB%AF-su-%C4%91%C3%BAc-cao-c%E1%BA%A5p-c%E1%BA%A1p-cao-t%C3%A0n-h%C3%ACnh-ch%E1%BB%AF-GODDESS-si%C3%AAu-%C4%91%E1%BA%B9p-i.237396802.22956820818?extraParams=%7B%22display_model_id%22%3A245036196414%2C%22model_selection_logic%22%3A3%7D&sp_atk=3c207934-60e7-4630-80a1-25caa0819eb0&xptdk=3c207934-60e7-4630-80a1-25caa0819eb0"

$ws.Hyperlinks.Add($ws.Range("A5"), $url, [Type]::Missing, [Type]::Missing, $url)

# Row heights for the header block (rows 1-6) grew slightly (15pt, thick
# bottom border) when the rows above the link were touched.
$ws.Range("A1:Z6").RowHeight = 15

# Selection moved to D14 in the saved file.
$ws.Range("D14").Select()
